# Fruta / hortaliza, semanal
# Insert two new weekly price rows (Primera / Segunda) at the top of the
# "Palta" data block, pushing the existing rows 640-669 down to 642-671.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at row 640 (shifts old rows 640.. down by 2)
$ws.Rows("640:641").Insert()

# --- New row 640: Primera ---
$ws.Cells.Item(640, 1).Value2 = 7
$ws.Cells.Item(640, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(640, 3).Value = "Ñuble"
$ws.Cells.Item(640, 4).Value2 = 44939
$ws.Cells.Item(640, 5).Value2 = 16
$ws.Cells.Item(640, 6).Value = "Fruta"
$ws.Cells.Item(640, 7).Value2 = 100106
$ws.Cells.Item(640, 8).Value = "Oleaginosos"
$ws.Cells.Item(640, 9).Value2 = 100106002
$ws.Cells.Item(640, 10).Value = "Palta"
$ws.Cells.Item(640, 11).Value = "Hass"
$ws.Cells.Item(640, 12).Value = "Primera"
$ws.Cells.Item(640, 13).Value2 = 200
$ws.Cells.Item(640, 14).Value2 = 3400
$ws.Cells.Item(640, 15).Value2 = 3500
$ws.Cells.Item(640, 16).Value2 = 3450
$ws.Cells.Item(640, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(640, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(640, 19).Value2 = 3450
$ws.Cells.Item(640, 20).Value2 = 1

# --- New row 641: Segunda ---
$ws.Cells.Item(641, 1).Value2 = 7
$ws.Cells.Item(641, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(641, 3).Value = "Ñuble"
$ws.Cells.Item(641, 4).Value2 = 44939
$ws.Cells.Item(641, 5).Value2 = 16
$ws.Cells.Item(641, 6).Value = "Fruta"
$ws.Cells.Item(641, 7).Value2 = 100106
$ws.Cells.Item(641, 8).Value = "Oleaginosos"
$ws.Cells.Item(641, 9).Value2 = 100106002
$ws.Cells.Item(641, 10).Value = "Palta"
$ws.Cells.Item(641, 11).Value = "Hass"
$ws.Cells.Item(641, 12).Value = "Segunda"
$ws.Cells.Item(641, 13).Value2 = 160
$ws.Cells.Item(641, 14).Value2 = 3100
$ws.Cells.Item(641, 15).Value2 = 3200
$ws.Cells.Item(641, 16).Value2 = 3150
$ws.Cells.Item(641, 17).Value = "`$/kilo (en caja de 15 kilos)"
$ws.Cells.Item(641, 18).Value = "Provincia de Quillota"
$ws.Cells.Item(641, 19).Value2 = 3150
$ws.Cells.Item(641, 20).Value2 = 1
